$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 used to only hold G21 = 1766.1 (a leftover balance value).
# Now it becomes a labeled test row: A21 = test name, H21 = "Test - value", I21 = amount.
$ws.Range("G21").ClearContents()
$ws.Range("A21").Value = "test_whenGetLastTransactions_thenGetLastNTransactionsAndAnteriorBalance"
$ws.Range("H21").Value = "Test – value"
$ws.Range("I21").Value = 1766.1

# Row 29 (H29 = 761.29) is moved up to row 22, with the same shape as row 21.
$ws.Range("A29:K29").EntireRow.Delete()
$ws.Range("A22").Value = "test_whenGetTransaction_thenTransactionHaveBalance"
$ws.Range("H22").Value = "Test – value"
$ws.Range("I22").Value = 761.29

# Move the active selection the way it ended up after the edit.
$ws.Range("J29").Select()
